$wb = $excel.ActiveWorkbook

# Insert the new "extreme_vols" worksheet between "vol_surface" and
# "bootstrapped_vol_surface", and give it the same look (styles/number
# formats/formulas) as "constant_vol_surface" by copying it.
$source = $wb.Worksheets.Item("constant_vol_surface")
$afterSheet = $wb.Worksheets.Item("vol_surface")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "extreme_vols"
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

$source.Range("A1:B13").Copy($newSheet.Range("A1"))

# Re-assert the tenor formulas in column A (Copy can drop these down to
# plain literals) so they stay live formulas, same as every other sheet.
$newSheet.Range("A3").Formula = "=1/12"
$newSheet.Range("A4").Formula = "=2/12"
$newSheet.Range("A5").Formula = "=3/12"
$newSheet.Range("A6").Formula = "=6/12"
$newSheet.Range("A7").Formula = "=9/12"

# Extreme original vol quotes for this new scenario.
$newSheet.Range("B2").Value = 10
$newSheet.Range("B3").Value = 12.333
$newSheet.Range("B4").Value = 14.154
$newSheet.Range("B5").Value = 15
$newSheet.Range("B6").Value = 20
$newSheet.Range("B7").Value = 25
$newSheet.Range("B8").Value = 30
$newSheet.Range("B9").Value = 35
$newSheet.Range("B10").Value = 40
$newSheet.Range("B11").Value = 50
$newSheet.Range("B12").Value = 70
$newSheet.Range("B13").Value = 80

# Mark the new sheet as the active/selected tab (it becomes the active
# tab of the workbook, matching the source edit).
$newSheet.Range("H9").Select()
$newSheet.Activate()
